# Apply the "product_stock_data" sheet update:
#  - decrement a few quantity values
#  - convert row 13's quantity/purchase_price/sale_price from text to real numbers
#  - append a new inventory row (14) for "Papel Contact Pliego"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- quantity corrections on existing rows ---
$ws.Range("D2").Value = 4
$ws.Range("D4").Value = 7
$ws.Range("D6").Value = 44
$ws.Range("D8").Value = 19
$ws.Range("D9").Value = 7

# --- row 13: make quantity/purchase_price/sale_price real numbers (were text) ---
$ws.Range("D13").Value = 12
$ws.Range("E13").Value = 1200
$ws.Range("F13").Value = 2400
$ws.Range("G13").Value = 45818.98260267361

# --- row 14: new product row ---
$ws.Range("A14").Value = "4P"
$ws.Range("B14").Value = "Papelería"
$ws.Range("C14").Value = "Papel Contact Pliego"

# D14/E14/F14 must stay text (like the legacy row 13 used to be) rather than
# becoming numeric, so force a text number format before assigning, then
# drop the format override back to the sheet default afterwards.
$ws.Range("D14:F14").NumberFormat = "@"
$ws.Range("D14").Value = "12"
$ws.Range("E14").Value = "400"
$ws.Range("F14").Value = "650"
$ws.Range("D14:F14").Style = "Normal"

$ws.Range("G14").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("G14").Value = 45821.65969966252

Write-Output "applied product_stock_data edits"
